# edit.ps1 - Applies the "Sprint review 1/2 done + udlæg til arbejdsplan" commit
# 1) Remove the stray _GoBack bookmark left after "Der skal afholdes et planning
#    meeting for sprint 1." (it is hidden from the normal Bookmarks enumeration,
#    but is still directly addressable by name).
$d = $word.ActiveDocument

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2) Replace the final (empty, ind-only) paragraph with the new "Mandag
#    10/12-2018" section: a Subtitle heading followed by the new bulleted
#    agenda items, ending with a fresh empty List Paragraph. InsertXML
#    replaces the contents of the range it is called on, so targeting the
#    last paragraph's Range swaps it out for the new paragraphs below
#    (one of which re-creates the _GoBack bookmark on "Opdater burndown-chart").
$last = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$last.InsertXML("<w:p><w:pPr><w:pStyle w:val=`"Undertitel`"/></w:pPr><w:r><w:t>Mandag 10/12-2018</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"Listeafsnit`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t xml:space=`"preserve`">Synkroniser </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>repositorys</w:t></w:r><w:proofErr w:type=`"spellEnd`"/></w:p><w:p><w:pPr><w:pStyle w:val=`"Listeafsnit`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t>Opdater burndown-chart</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p><w:p><w:pPr><w:pStyle w:val=`"Listeafsnit`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t xml:space=`"preserve`">Rollefordelingsdiskussion: Patrick som product </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>owner</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`">, X som stand-in product </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>owner</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> under review og Patrick som stakeholder.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"Listeafsnit`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:proofErr w:type=`"spellStart`"/><w:r><w:t>Acceptance</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t>-test vs. Hvis/Når/Så-test: formål og go/no go.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"Listeafsnit`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t>Færdiggøre afsnit om sprint review og retrospekt (</w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>casper</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t>)</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"Listeafsnit`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t>Gennemgang på dagen hvis færdiggjort.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"Listeafsnit`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:proofErr w:type=`"spellStart`"/><w:r><w:t>Reestimering</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> (tilføj nye kolonner til alle</w:t></w:r><w:r><w:t xml:space=`"preserve`"> ikke-afsluttede</w:t></w:r><w:r><w:t xml:space=`"preserve`"> User Storys som viser 2. estimering).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"Listeafsnit`"/></w:pPr></w:p>")
